$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("E3").Value = "2016-03-22 08:43:05"
$zh.Range("H3").Value = "2016-03-22 08:43:29"

$de = $wb.Worksheets.Item("de-de")
$de.Range("E3").Value = "2016-03-22 08:43:09"
$de.Range("H3").Value = "2016-03-22 08:43:36"
